$wb = $excel.ActiveWorkbook

# ----- Sheet 1: LP1912 -----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 03:42:43"
$ws1.Range("A3").Value = "Total filas: 18"

$ws1.Cells.Item(6, 1).Value = "01:12:01"
$ws1.Cells.Item(6, 2).Value = "01:12"
$ws1.Cells.Item(6, 3).Value = "215_ALUAR"
$ws1.Cells.Item(6, 4).Value = 0
$ws1.Cells.Item(6, 5).Value = "LP1912"

$ws1.Cells.Item(7, 1).Value = "01:56:31"
$ws1.Cells.Item(7, 2).Value = "01:58"
$ws1.Cells.Item(7, 3).Value = "14_ABASTO"
$ws1.Cells.Item(7, 4).Value = 2
$ws1.Cells.Item(7, 5).Value = "LP1912"

$ws1.Cells.Item(8, 1).Value = "02:49:45"
$ws1.Cells.Item(8, 2).Value = "02:49"
$ws1.Cells.Item(8, 3).Value = "215_ALUAR"
$ws1.Cells.Item(8, 4).Value = 0
$ws1.Cells.Item(8, 5).Value = "LP1912"

$ws1.Cells.Item(9, 1).Value = "01:12:01"
$ws1.Cells.Item(9, 2).Value = "02:58"
$ws1.Cells.Item(9, 3).Value = "215_ALUAR"
$ws1.Cells.Item(9, 4).Value = 106
$ws1.Cells.Item(9, 5).Value = "LP1912"

$ws1.Cells.Item(10, 1).Value = "01:56:31"
$ws1.Cells.Item(10, 2).Value = "02:59"
$ws1.Cells.Item(10, 3).Value = "215_ALUAR"
$ws1.Cells.Item(10, 4).Value = 63
$ws1.Cells.Item(10, 5).Value = "LP1912"

$ws1.Cells.Item(11, 1).Value = "03:42:43"
$ws1.Cells.Item(11, 2).Value = "03:45"
$ws1.Cells.Item(11, 3).Value = "14_ABASTO"
$ws1.Cells.Item(11, 4).Value = 3
$ws1.Cells.Item(11, 5).Value = "LP1912"

$ws1.Cells.Item(12, 1).Value = "03:00:18"
$ws1.Cells.Item(12, 2).Value = "03:48"
$ws1.Cells.Item(12, 3).Value = "14_ABASTO"
$ws1.Cells.Item(12, 4).Value = 48
$ws1.Cells.Item(12, 5).Value = "LP1912"

$ws1.Cells.Item(13, 1).Value = "02:24:16"
$ws1.Cells.Item(13, 2).Value = "03:53"
$ws1.Cells.Item(13, 3).Value = "14_ABASTO"
$ws1.Cells.Item(13, 4).Value = 89
$ws1.Cells.Item(13, 5).Value = "LP1912"

$ws1.Cells.Item(14, 1).Value = "02:24:16"
$ws1.Cells.Item(14, 2).Value = "03:58"
$ws1.Cells.Item(14, 3).Value = "215_ALUAR"
$ws1.Cells.Item(14, 4).Value = 94
$ws1.Cells.Item(14, 5).Value = "LP1912"

$ws1.Cells.Item(15, 1).Value = "03:42:43"
$ws1.Cells.Item(15, 2).Value = "04:01"
$ws1.Cells.Item(15, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(15, 4).Value = 19
$ws1.Cells.Item(15, 5).Value = "LP1912"

$ws1.Cells.Item(16, 1).Value = "02:49:45"
$ws1.Cells.Item(16, 2).Value = "04:35"
$ws1.Cells.Item(16, 3).Value = "215_ALUAR"
$ws1.Cells.Item(16, 4).Value = 106
$ws1.Cells.Item(16, 5).Value = "LP1912"

$ws1.Cells.Item(17, 1).Value = "03:00:18"
$ws1.Cells.Item(17, 2).Value = "04:44"
$ws1.Cells.Item(17, 3).Value = "215_ALUAR"
$ws1.Cells.Item(17, 4).Value = 104
$ws1.Cells.Item(17, 5).Value = "LP1912"

$ws1.Cells.Item(18, 1).Value = "03:42:43"
$ws1.Cells.Item(18, 2).Value = "04:45"
$ws1.Cells.Item(18, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(18, 4).Value = 63
$ws1.Cells.Item(18, 5).Value = "LP1912"

$ws1.Cells.Item(19, 1).Value = "03:42:43"
$ws1.Cells.Item(19, 2).Value = "04:53"
$ws1.Cells.Item(19, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(19, 4).Value = 71
$ws1.Cells.Item(19, 5).Value = "LP1912"

$ws1.Cells.Item(20, 1).Value = "03:42:43"
$ws1.Cells.Item(20, 2).Value = "05:16"
$ws1.Cells.Item(20, 3).Value = "17_ROMERO"
$ws1.Cells.Item(20, 4).Value = 94
$ws1.Cells.Item(20, 5).Value = "LP1912"

$ws1.Cells.Item(21, 1).Value = "03:42:43"
$ws1.Cells.Item(21, 2).Value = "05:22"
$ws1.Cells.Item(21, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(21, 4).Value = 100
$ws1.Cells.Item(21, 5).Value = "LP1912"

$ws1.Cells.Item(22, 1).Value = "03:42:43"
$ws1.Cells.Item(22, 2).Value = "05:34"
$ws1.Cells.Item(22, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(22, 4).Value = 112
$ws1.Cells.Item(22, 5).Value = "LP1912"

$ws1.Cells.Item(23, 1).Value = "03:42:43"
$ws1.Cells.Item(23, 2).Value = "05:35"
$ws1.Cells.Item(23, 3).Value = "14_ABASTO"
$ws1.Cells.Item(23, 4).Value = 113
$ws1.Cells.Item(23, 5).Value = "LP1912"

# ----- Sheet 2: LP1912-215 -----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 03:42:43"
$ws2.Range("A3").Value = "Total filas: 9"

$ws2.Cells.Item(13, 1).Value = "03:42:43"
$ws2.Cells.Item(13, 2).Value = "04:45"
$ws2.Cells.Item(13, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(13, 4).Value = 63
$ws2.Cells.Item(13, 5).Value = "LP1912"

$ws2.Cells.Item(14, 1).Value = "03:42:43"
$ws2.Cells.Item(14, 2).Value = "05:34"
$ws2.Cells.Item(14, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(14, 4).Value = 112
$ws2.Cells.Item(14, 5).Value = "LP1912"

# ----- Sheet 3: 6203-6173 -----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 03:42:43"

